# Auto-generated edit script derived from the target diff for Saldo.xlsx.
#
# The sheet is a flat account/name/balance table with no explicit row
# numbers (implicit top-to-bottom layout). The diff is a series of
# localized row block replacements/inserts/deletes scattered through the
# table, so we walk them top-to-bottom (in ascending row order) and, for
# each block:
#   - insert/delete whole rows to make room / remove rows as needed
#   - (re)write the Conta/Nome/Saldo values for every row in the block
#
# Column A ('Conta') holds zero-padded account numbers (e.g. 005135105);
# without help Excel's Value setter treats them as numbers and drops the
# leading zeros, so we write them with a leading apostrophe to force text,
# then ClearFormats() to drop the 'Number Stored as Text' quote-prefix
# style Excel attaches, matching the plain (unstyled) text cells used for
# this column everywhere else in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# replace: rows 2.. (old_count=1, new_count=1)
$ws.Cells.Item(2,1).Value = "'005135105"
$ws.Cells.Item(2,1).ClearFormats()
$ws.Cells.Item(2,2).Value = "Brenner"
$ws.Cells.Item(2,3).Value = 208950.05

# replace: rows 4.. (old_count=3, new_count=4)
$ws.Range("4:4").Insert()
$ws.Cells.Item(4,1).Value = "'004935287"
$ws.Cells.Item(4,1).ClearFormats()
$ws.Cells.Item(4,2).Value = "Odilon"
$ws.Cells.Item(4,3).Value = 58025.61
$ws.Cells.Item(5,1).Value = "'000806386"
$ws.Cells.Item(5,1).ClearFormats()
$ws.Cells.Item(5,2).Value = "Fernanda"
$ws.Cells.Item(5,3).Value = 50533.82
$ws.Cells.Item(6,1).Value = "'004444380"
$ws.Cells.Item(6,1).ClearFormats()
$ws.Cells.Item(6,2).Value = "Marcelo"
$ws.Cells.Item(6,3).Value = 49648.92
$ws.Cells.Item(7,1).Value = "'005599726"
$ws.Cells.Item(7,1).ClearFormats()
$ws.Cells.Item(7,2).Value = "Jorge"
$ws.Cells.Item(7,3).Value = 26275.75

# delete: rows 9.. (old_count=1, new_count=0)
$ws.Range("9:9").Delete()

# delete: rows 10.. (old_count=4, new_count=0)
$ws.Range("10:13").Delete()

# replace: rows 11.. (old_count=1, new_count=2)
$ws.Range("11:11").Insert()
$ws.Cells.Item(11,1).Value = "'004755083"
$ws.Cells.Item(11,1).ClearFormats()
$ws.Cells.Item(11,2).Value = "Evaldo"
$ws.Cells.Item(11,3).Value = 5994.58
$ws.Cells.Item(12,1).Value = "'004202332"
$ws.Cells.Item(12,1).ClearFormats()
$ws.Cells.Item(12,2).Value = "Tatiana"
$ws.Cells.Item(12,3).Value = 4985.28

# insert: rows 18.. (old_count=0, new_count=2)
$ws.Range("18:19").Insert()
$ws.Cells.Item(18,1).Value = "'008404765"
$ws.Cells.Item(18,1).ClearFormats()
$ws.Cells.Item(18,2).Value = "Leticia"
$ws.Cells.Item(18,3).Value = 1000
$ws.Cells.Item(19,1).Value = "'004355790"
$ws.Cells.Item(19,1).ClearFormats()
$ws.Cells.Item(19,2).Value = "Mineia"
$ws.Cells.Item(19,3).Value = 976.92

# insert: rows 28.. (old_count=0, new_count=1)
$ws.Range("28:28").Insert()
$ws.Cells.Item(28,1).Value = "'005198093"
$ws.Cells.Item(28,1).ClearFormats()
$ws.Cells.Item(28,2).Value = "Ana"
$ws.Cells.Item(28,3).Value = 499.1

# delete: rows 42.. (old_count=1, new_count=0)
$ws.Range("42:42").Delete()
